$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.197.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.242.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.79%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -15.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.94%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.240.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.61%  "

$ws.Range("E10").Value = "  -9.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.130"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.59%  "

$ws.Range("E14").Value = "  -10.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.752.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.238.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.62%  "

$ws.Range("E17").Value = "  -8.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.061.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.939"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "364.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.63%  "

$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.05%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "626.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.63%  "

$ws.Range("E34").Value = "  -6.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.29%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0669"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.865.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.120"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.31%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.32%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.37%  "

$ws.Range("E46").Value = "  -14.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0383"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.122"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.36%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.44%  "
